{"js": "// Office.js (Word JavaScript API) script.\n// Adds the \"handCount\" formula write-up right after the paragraph that\n// ends with \"...so I couldn't get it exactly how I wanted it.\" and\n// before the trailing (_GoBack-bookmarked) blank paragraph at the end\n// of the document.\n\n// 1. Locate the target paragraph by its distinctive trailing text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"so I couldn\\u2019t get it exactly how I wanted it.\") !== -1) {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the anchor paragraph for the insertion.\");\n}\nconst target = paragraphs.items[targetIndex];\n\n// 2. The document carries a hidden \"_GoBack\" bookmark inside that same\n//    paragraph. The new content needs to end up between the paragraph's\n//    text and that bookmark, so pull the bookmark out first and we will\n//    re-create it afterwards, in its new home.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3. Create an anchor paragraph right after the target paragraph; we will\n//    replace its contents (via OOXML) with the whole block of new\n//    paragraphs so the run/tab/bookmark structure comes out exactly right.\nconst anchor = target.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nconst anchorRange = anchor.getRange(Word.RangeLocation.content);\n\nconst newContentOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">I made a function called </w:t></w:r>' +\n  '<w:r><w:rPr><w:i/></w:rPr><w:t>handCount</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t>function handCount(index) {</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:tab/><w:t>for(var I = 0; I &lt; index; i++) {</w:t></w:r></w:p>' +\n  '<w:p>' +\n  '<w:r><w:tab/></w:r>' +\n  '<w:r><w:tab/><w:t>if ((I \\u2013 1) % 8 == 0) {</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:r><w:tab/></w:r>' +\n  '<w:r><w:tab/></w:r>' +\n  '<w:r><w:tab/><w:t>console.log(i)</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\"> }}}</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nanchorRange.insertOoxml(newContentOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Adds the \"handCount\" formula write-up right after the paragraph that\n# ends with \"...so I couldn't get it exactly how I wanted it.\" and\n# before the trailing (_GoBack-bookmarked) blank paragraph at the end\n# of the document.\n\n$d = $word.ActiveDocument\n\n# 1. Locate the target paragraph by its distinctive trailing text.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*so I couldn*t get it exactly how I wanted it.*\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not find the anchor paragraph for the insertion.\"\n}\n$target = $d.Paragraphs.Item($targetIndex)\n\n# 2. The document carries a hidden \"_GoBack\" bookmark inside that same\n#    paragraph. The new content needs to end up between the paragraph's\n#    text and that bookmark, so remove the bookmark first; it gets\n#    re-created afterwards, in its new home.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 3. Create an anchor paragraph right after the target paragraph, then\n#    replace its contents (via WordOpenXML/InsertXML) with the whole\n#    block of new paragraphs so the run/tab/bookmark structure comes out\n#    exactly right.\n$r = $target.Range\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n$anchor = $d.Paragraphs.Item($targetIndex + 1)\n$anchorRange = $anchor.Range\n\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">I made a function called </w:t></w:r>\n            <w:r><w:rPr><w:i/></w:rPr><w:t>handCount</w:t></w:r>\n          </w:p>\n          <w:p/>\n          <w:p/>\n          <w:p/>\n          <w:p/>\n          <w:p/>\n          <w:p/>\n          <w:p>\n            <w:r><w:t>function handCount(index) {</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:tab/><w:t>for(var I = 0; I &lt; index; i++) {</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:tab/></w:r>\n            <w:r><w:tab/><w:t>if ((I &#8211; 1) % 8 == 0) {</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:tab/></w:r>\n            <w:r><w:tab/></w:r>\n            <w:r><w:tab/><w:t>console.log(i)</w:t></w:r>\n          </w:p>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\"> }}}</w:t></w:r>\n          </w:p>\n          <w:p/>\n          <w:p>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$anchorRange.InsertXML($ooxml)\n"}
